$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Status" column: header + first couple of data values
$ws.Range("D1").Value = "Status"
$ws.Range("D2").Value = "Done"
$ws.Range("D3").Value = "Working"

# D1 is a header cell like A1/B1/C1 -> bold font
$ws.Range("D1").Font.Bold = $true

# Header row (A1:D1) gets a thin black bottom border
$b = $ws.Range("A1:D1").Borders.Item(9)
$b.LineStyle = 1
$b.Weight = 2
$b.Color = 0
